{"js": "// Insert \"-C\" into \"(Figure 2B)\" to make it \"(Figure 2B-C)\", matching the\n// author's edit that added a reference to an additional sub-panel (Figure 2C)\n// of the normalized AUC statistic. The edit point also carries forward\n// Word's \"last edit\" (_GoBack) bookmark, which the author's edit left sitting\n// immediately after the newly typed \"-C\" (i.e. right before the closing\n// parenthesis).\n\n// 1) Drop the existing _GoBack bookmark -- Word always re-anchors this\n//    bookmark at the most recent edit location, so the one currently in the\n//    document (at the very end of the sentence) needs to move.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the unique occurrence of \"Figure 2B\" (without the parenthesis) so\n//    we can collapse a range right after the \"B\" and before the \")\".\nconst searchResults = context.document.body.search(\"Figure 2B\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length !== 1) {\n  throw new Error(`Expected exactly one match for \"Figure 2B\", found ${searchResults.items.length}`);\n}\n\nconst insertionPoint = searchResults.items[0].getRange(\"End\");\n\n// 3) Type the new text at that location.\ninsertionPoint.insertText(\"-C\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Re-create the _GoBack bookmark immediately after what was just typed\n//    (i.e. right before the closing parenthesis), matching where Word leaves\n//    the mark after an in-place text edit.\nconst afterInsert = context.document.body.search(\"Figure 2B-C\", { matchCase: true });\nafterInsert.load(\"text\");\nawait context.sync();\n\nif (afterInsert.items.length !== 1) {\n  throw new Error(`Expected exactly one match for \"Figure 2B-C\", found ${afterInsert.items.length}`);\n}\n\nafterInsert.items[0].getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Insert \"-C\" into \"(Figure 2B)\" to make it \"(Figure 2B-C)\", matching the\n# author's edit that added a reference to an additional sub-panel (Figure 2C)\n# of the normalized AUC statistic. The edit point also carries forward\n# Word's \"last edit\" (_GoBack) bookmark, which the author's edit left sitting\n# immediately after the newly typed \"-C\" (i.e. right before the closing\n# parenthesis).\n\n$d = $word.ActiveDocument\n\n# 1) Drop the existing _GoBack bookmark -- Word always re-anchors this\n#    bookmark at the most recent edit location, so the one currently in the\n#    document (at the very end of the sentence) needs to move.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Locate the unique occurrence of \"Figure 2B\" (without the parenthesis) so\n#    we can collapse a range right after the \"B\" and before the \")\".\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Figure 2B\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find 'Figure 2B'\"\n}\n\n$insertStart = $rng.End\n$rng.Collapse(0)  # wdCollapseEnd\n\n# 3) Type the new text at that location.\n$rng.InsertAfter(\"-C\")\n\n# 4) Re-create the _GoBack bookmark immediately after what was just typed\n#    (i.e. right before the closing parenthesis), matching where Word leaves\n#    the mark after an in-place text edit. Built from absolute offsets so it\n#    exactly spans the two characters we just inserted (\"-C\").\n$newBookmarkPoint = $d.Range($insertStart + 2, $insertStart + 2)\n$d.Bookmarks.Add(\"_GoBack\", $newBookmarkPoint)\n"}
